# Apply the threshold-table update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data values (rows 2, 3 and 5) ---
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 11.5
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 10.5
$ws.Range("C5").Value = 20

# --- Add the new (empty but styled) block G7:H10 ---
# Copy an existing, already-styled cell and paste it into the new block so
# the new cells reuse the same cell style (s="2") instead of Excel creating
# a brand-new style entry, then clear the copied content so the cells stay
# empty just like in the target sheet.
$ws.Range("A2").Copy()
$ws.Range("G7:H10").PasteSpecial(-4122)
$ws.Range("G7:H10").ClearContents()
$excel.CutCopyMode = $false | Out-Null

# --- Update the selected/active cell shown when the sheet is opened ---
$ws.Range("E9").Select() | Out-Null
